$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Rushing  (Week 16 simulation updates)
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# D.Carr (row 2)
$rushing.Range("C2").Value = 10
$rushing.Range("E2").Value = 10
$rushing.Range("F2").Value = 6

# M.Mariota (row 3)
$rushing.Range("D3").Value = 1

# J.Jacobs (row 5)
$rushing.Range("C5").Value = 75
$rushing.Range("D5").Value = 55
$rushing.Range("E5").Value = 8

# P.Barber (row 7)
$rushing.Range("C7").Value = 32
$rushing.Range("D7").Value = 15

# ---------------------------------------------------------------
# Sheet: Receiving (Week 16 simulation updates + Week 15 log)
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# J.Jacobs (row 2)
$receiving.Range("C2").Value = 51
$receiving.Range("D2").Value = 44

# A.Ingold (row 4)
$receiving.Range("C4").Value = 12
$receiving.Range("D4").Value = 10
$receiving.Range("G4").Value = 2
$receiving.Range("H4").Value = 2

# B.Edwards (row 7)
$receiving.Range("C7").Value = 26
$receiving.Range("D7").Value = 18
$receiving.Range("G7").Value = 6
$receiving.Range("H7").Value = 3

# H.Renfrow (row 8)
$receiving.Range("C8").Value = 97
$receiving.Range("D8").Value = 79
$receiving.Range("G8").Value = 18
$receiving.Range("H8").Value = 12

# Z.Jones (row 9)
$receiving.Range("C9").Value = 22
$receiving.Range("D9").Value = 18
$receiving.Range("E9").Value = 18

# D.Jackson (row 11)
$receiving.Range("C11").Value = 5
$receiving.Range("D11").Value = 4
$receiving.Range("E11").Value = 6

# F.Moreau (row 13)
$receiving.Range("C13").Value = 31
$receiving.Range("D13").Value = 20

# New player row 15: D.Helm
# Copy formatting from the row above (bold/centered/bordered "index" style)
# before writing values, so the new row matches the existing table styling.
$receiving.Range("A14").Copy()
$receiving.Range("A15").PasteSpecial(-4122)
$receiving.Range("A15").Value = 13
$receiving.Range("B15").Value = "D.Helm"
$receiving.Range("C15").Value = 1
$receiving.Range("D15").Value = 1
$receiving.Range("E15").Value = 0
$receiving.Range("F15").Value = 0
$receiving.Range("G15").Value = 0
$receiving.Range("H15").Value = 0
